$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 27
$ws1.Range("F4").Value = 19
$ws1.Range("F5").Value = 15870
$ws1.Range("F6").Value = 421
$ws1.Range("F7").Value = 11
$ws1.Range("F9").Value = 15481
$ws1.Range("F11").Value = 9097
$ws1.Range("F12").Value = 395
$ws1.Range("F14").Value = 1019
$ws1.Range("F15").Value = 107
$ws1.Range("F16").Value = 206
$ws1.Range("F18").Value = 207
$ws1.Range("F20").Value = 64
$ws1.Range("F21").Value = 569
$ws1.Range("F22").Value = 27
$ws1.Range("F25").Value = 1117
$ws1.Range("F30").Value = 34
$ws1.Range("F36").Value = 330
$ws1.Range("F39").Value = 5589

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 27
$ws4.Range("F4").Value = 19
$ws4.Range("F5").Value = 15870
$ws4.Range("F6").Value = 421
$ws4.Range("F7").Value = 11
$ws4.Range("F9").Value = 15481
$ws4.Range("F11").Value = 9097
$ws4.Range("F12").Value = 395
$ws4.Range("F14").Value = 1019
$ws4.Range("F15").Value = 107
$ws4.Range("F16").Value = 206
$ws4.Range("F18").Value = 207
$ws4.Range("F20").Value = 64
$ws4.Range("F21").Value = 569
$ws4.Range("F22").Value = 27
$ws4.Range("F25").Value = 1117
$ws4.Range("F30").Value = 34
$ws4.Range("F38").Value = 330
$ws4.Range("F41").Value = 5589
